$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows with an ambiguous numeric-looking Price (column D) value need a leading
# apostrophe so Excel keeps them as Text (matching the source's inlineStr cells)
# instead of auto-converting to a Number; resetting Style to Normal afterward
# clears the quote-prefix style flag Excel would otherwise stamp on the cell.

$ws.Range("D2").Value = "67.166.01"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "3.518.25"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'593.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "'173.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.595"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.46%  "
$ws.Range("E9").Value = "  +6.84%  "
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "4.125.04"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "'29.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").Value = "67.124.21"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "3.478.03"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "'6.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "'14.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").Value = "'395.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "'0.541"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("E25").Value = "  -3.90%  "
$ws.Range("D26").Value = "'10.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.70%  "
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "'6.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").Value = "'2.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").Value = "'23.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("D33").Value = "'7.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").Value = "'1.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("D35").Value = "'162.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").Value = "'6.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.38%  "
$ws.Range("D39").Value = "'27.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.94%  "
$ws.Range("D40").Value = "'4.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").Value = "'26.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").Value = "'2.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.36%  "
$ws.Range("D44").Value = "2.805.33"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("D47").Value = "'336.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.15%  "
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").Value = "'33.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "'0.850"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.35%  "
